# mod: fit symbol to firstline length
#
# Rewrites the small Sheet1 "No / sample1 / long name test pattern" table:
#   - header row keeps its header style, but becomes text labels
#   - data rows lose the special "thick divider" row (old row3) / "totals"
#     row (old row4) styling and all collapse onto the same plain bordered
#     style that row 2 already used
#   - row heights go back to the sheet's normal (default) height
#   - selection moves to C6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. New cell content -----------------------------------------------
# Row 1: header labels (kept centered/shaded "header" style s=2)
$ws.Range("A1").Value2 = "No"
$ws.Range("B1").Value2 = "sample1"
$ws.Range("C1").Value2 = "long name test pattern"

# Row 2
$ws.Range("A2").Value2 = 1
$ws.Range("B2").Value2 = "test1"
$ws.Range("C2").Value2 = "long long test"

# Row 3 (previously the bold/thick-bottom divider row)
$ws.Range("A3").Value2 = 2
$ws.Range("B3").Value2 = "test2"
$ws.Range("C3").Value2 = "long long test"

# Row 4 (previously the numeric "totals" row)
$ws.Range("A4").Value2 = 3
$ws.Range("B4").Value2 = "test3"
$ws.Range("C4").Value2 = "long long test"

# --- 2. Normalise formatting of rows 3 & 4 to match row 2 ---------------
# Use copy/paste-special (formats only) instead of touching borders
# directly so the existing "plain bordered" style already used by row 2
# is reused verbatim rather than a near-duplicate style being created.
$ws.Range("A2").Copy()
$ws.Range("A3:C4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Drop the custom row heights / thick top+bottom borders ----------
# Rows 3 & 4 previously had an explicit 19.5pt height (with thick
# bottom/top borders); put them back to the sheet's normal row height.
$normalHeight = $ws.Rows.Item(1).RowHeight
$ws.Rows.Item(3).RowHeight = $normalHeight
$ws.Rows.Item(3).EntireRow.AutoFit()
$ws.Rows.Item(4).RowHeight = $normalHeight
$ws.Rows.Item(4).EntireRow.AutoFit()

# --- 4. Update the active selection -------------------------------------
$null = $ws.Range("C6").Select()
